# Applies the two reachable effects of the authored change:
#   1. The table on slide 16 switches from the deck's default/custom
#      table style to the built-in style {2F221D35-A5D8-4A76-9F06-51A10F55E3FC}.
#   2. The presentation's theme (ppt/theme/theme1.xml, the theme used by the
#      slide master / every slide) switches its 12-slot colour scheme from the
#      "Integral" palette to the stock "Office" palette - i.e. the deck's
#      Design colours become the default Office theme colours.

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # VBA/PowerPoint RGB() long: 0x00BBGGRR
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{2F221D35-A5D8-4A76-9F06-51A10F55E3FC}")
    }
}

# --- 2. Theme colour scheme -> stock "Office" colours ----------------------
# Order matches a:clrScheme child order / ThemeColorScheme 1-based index:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Item($i).RGB = HexToRgbInt $officeColors[$i - 1]
}
